# Apply updated odds values to Sheet1, matching the target diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 3
$ws.Range("M3").Value = 1.13
$ws.Range("O3").Value = 1.62
$ws.Range("P3").Value = 2.3

# Row 5
$ws.Range("J5").Value = 3.6
$ws.Range("Q5").Value = 2
$ws.Range("R5").Value = 1.85
$ws.Range("AA5").Value = 23
$ws.Range("AH5").Value = 8
$ws.Range("AJ5").Value = 9.5

# Row 8
$ws.Range("G8").Value = 2.4
$ws.Range("H8").Value = 3
$ws.Range("I8").Value = 3
$ws.Range("J8").Value = 3.1
$ws.Range("L8").Value = 3.5
$ws.Range("N8").Value = 9.5
$ws.Range("W8").Value = 8.5
$ws.Range("X8").Value = 12
$ws.Range("Y8").Value = 10
$ws.Range("Z8").Value = 23
$ws.Range("AA8").Value = 19
$ws.Range("AB8").Value = 29
$ws.Range("AE8").Value = 12
$ws.Range("AH8").Value = 10
$ws.Range("AI8").Value = 15
$ws.Range("AK8").Value = 29
$ws.Range("AL8").Value = 23
$ws.Range("AN8").Value = 4.5
$ws.Range("AO8").Value = 13
$ws.Range("AP8").Value = 23
$ws.Range("AR8").Value = 67
$ws.Range("AX8").Value = 17
$ws.Range("BA8").Value = 67
$ws.Range("BB8").Value = 151

# Row 9
$ws.Range("G9").Value = 1.44
$ws.Range("J9").Value = 1.91
$ws.Range("K9").Value = 2.63
$ws.Range("Q9").Value = 1.48
$ws.Range("R9").Value = 2.6
$ws.Range("U9").Value = 1.62
$ws.Range("V9").Value = 2.2
$ws.Range("W9").Value = 10
$ws.Range("AI9").Value = 41
$ws.Range("AR9").Value = 34
$ws.Range("BC9").Value = 126

# Row 10
$ws.Range("O10").Value = 1.44
$ws.Range("P10").Value = 2.63
$ws.Range("Q10").Value = 2.35
$ws.Range("R10").Value = 1.57

# Row 12
$ws.Range("O12").Value = 1.5
$ws.Range("P12").Value = 2.5

# Row 13
$ws.Range("M13").Value = 1.03
$ws.Range("N13").Value = 15
